$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

$ws.Range("F2").Value = 34.05887999999999
$ws.Range("K2").Value = 0

$ws.Range("J3").Value = 177.673824
$ws.Range("K3").Value = 1281.181536
$ws.Range("L3").Value = 6185.342524608289

$ws.Range("J4").Value = 34.05888
$ws.Range("K4").Value = 34.05888
$ws.Range("L4").Value = 1082.329036961193
